$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values regenerated to filter save games (row -> B,C,D,E,G)
$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 28.30127388105354)
    3 = @(0.3464964993005633, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 1.051601690082842)
    4 = @(0.06328177979961902, 0.05231270169004087, 3.082599426703578, 6.48142807727062, 9.679621985463859)
    5 = @(0.7287194209349384, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 19.86557370323023)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $vals[4]
}
